$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0) Make room: the table grows from 6 data rows (rows 2-7) to 8 data rows
#    (rows 2-9), so insert two blank rows at the bottom of the existing table.
# ---------------------------------------------------------------------------
$ws.Rows.Item("8:9").Insert()

# ---------------------------------------------------------------------------
# 1) Capture the three distinct cell-format "patterns" used in the table as
#    templates in a scratch area, before we start overwriting rows 2-7.
#      - row 2 pattern : F/G filled in,  column A normal style
#      - row 3 pattern : F/G blank,      column A normal style
#      - row 7 pattern : F/G blank,      column A "last of group" style
# ---------------------------------------------------------------------------
$ws.Range("A2:H2").Copy()
$ws.Range("A50:H50").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A3:H3").Copy()
$ws.Range("A51:H51").PasteSpecial(-4122)

$ws.Range("A7:H7").Copy()
$ws.Range("A52:H52").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Apply the right format template to every data row of the final table.
# ---------------------------------------------------------------------------
$filled = $ws.Range("A50:H50")
$blank  = $ws.Range("A51:H51")
$lastBlank = $ws.Range("A52:H52")

$filled.Copy();    $ws.Range("A2:H2").PasteSpecial(-4122)
$filled.Copy();    $ws.Range("A3:H3").PasteSpecial(-4122)
$filled.Copy();    $ws.Range("A4:H4").PasteSpecial(-4122)
$filled.Copy();    $ws.Range("A5:H5").PasteSpecial(-4122)
$filled.Copy();    $ws.Range("A6:H6").PasteSpecial(-4122)
$blank.Copy();     $ws.Range("A7:H7").PasteSpecial(-4122)
$blank.Copy();     $ws.Range("A8:H8").PasteSpecial(-4122)
$lastBlank.Copy(); $ws.Range("A9:H9").PasteSpecial(-4122)

# row 6's ItemKey cell uses the "last of group" style, same as the one used
# for rows 9 -- fix column A specifically.
$ws.Range("A52").Copy()
$ws.Range("A6").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Remove the scratch templates.
# ---------------------------------------------------------------------------
$ws.Range("A50:H52").Clear()

# ---------------------------------------------------------------------------
# 4) Row heights -- only rows that need to be taller than the sheet default
#    get an explicit height; the rest keep using the sheet's default height.
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(8).RowHeight = 30
$ws.Rows.Item(9).RowHeight = 30

# ---------------------------------------------------------------------------
# 4b) Make sure brand-new shared strings get allocated in the same order as
#     the reference workbook: SILT, Sandy Silt, Silty Clay.
# ---------------------------------------------------------------------------
$ws.Range("B3").Value2 = "SILT"
$ws.Range("A6").Value2 = "Sandy Silt"
$ws.Range("A3").Value2 = "Silty Clay"

# ---------------------------------------------------------------------------
# 5) Values / formulas for the (re-ordered & extended) table.
# ---------------------------------------------------------------------------
# Row 2 : Clay / CLAY (unchanged)
$ws.Range("A2").Value2 = "Clay"
$ws.Range("B2").Value2 = "CLAY"
$ws.Range("C2").Value2 = 20
$ws.Range("D2").Value2 = 0.5
$ws.Range("F2").Value2 = 15
$ws.Range("G2").Value2 = 20
$ws.Range("H2").Value2 = 215

# Row 3 : Silty Clay / SILT (plain value, no formula)
$ws.Range("A3").Value2 = "Silty Clay"
$ws.Range("B3").Value2 = "SILT"
$ws.Range("C3").Value2 = 20
$ws.Range("D3").Value2 = 0.5
$ws.Range("F3").Value2 = 15
$ws.Range("G3").Value2 = 20
$ws.Range("H3").Value2 = 215

# Row 4 : Sandy Clay / SILT (plain, non-shared, formula)
$ws.Range("A4").Value2 = "Sandy Clay"
$ws.Range("B4").Value2 = "SILT"
$ws.Range("C4").Formula = "=C2"
$ws.Range("D4").Value2 = 0.5
$ws.Range("F4").Value2 = 15
$ws.Range("G4").Value2 = 20
$ws.Range("H4").Value2 = 215

# Row 5 : Silt / SILT
$ws.Range("A5").Value2 = "Silt"
$ws.Range("B5").Value2 = "SILT"
$ws.Range("D5").Value2 = 0.5
$ws.Range("F5").Value2 = 15
$ws.Range("G5").Value2 = 20
$ws.Range("H5").Value2 = 215

# Row 6 : Sandy Silt / SILT
$ws.Range("A6").Value2 = "Sandy Silt"
$ws.Range("B6").Value2 = "SILT"
$ws.Range("D6").Value2 = 0.5
$ws.Range("F6").Value2 = 15
$ws.Range("G6").Value2 = 20
$ws.Range("H6").Value2 = 215

# Row 7 : Sand / SAND
$ws.Range("A7").Value2 = "Sand"
$ws.Range("B7").Value2 = "SAND"
$ws.Range("D7").Value2 = 0.5
$ws.Range("H7").Value2 = 215

# Row 8 : Silty Sand / SILT
$ws.Range("A8").Value2 = "Silty Sand"
$ws.Range("B8").Value2 = "SILT"
$ws.Range("D8").Value2 = 0.5
$ws.Range("H8").Value2 = 215

# Row 9 : Clayey Sand / SILT
$ws.Range("A9").Value2 = "Clayey Sand"
$ws.Range("B9").Value2 = "SILT"
$ws.Range("D9").Value2 = 0.5
$ws.Range("H9").Value2 = 215

# C5:C9 form a single shared formula block ( =C4, =C5, =C6, =C7, =C8 ),
# created in one shot so Excel records it as one shared-formula group.
$ws.Range("C5:C9").Formula = "=C4"

$excel.Calculate()

# ---------------------------------------------------------------------------
# 6) View bits.
# ---------------------------------------------------------------------------
$ws.Range("P6").Select()

$wb.Save()
